{"js": "// Append the \"Acknowledgments and Dedication\" heading as a new run in the\n// last (empty) paragraph of the document body, matching the author's\n// existing Arial/bold/14pt styling used throughout the title page.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst range = lastParagraph.insertText(\"Acknowledgments and Dedication\", \"End\");\nrange.font.set({\n  name: \"Arial\",\n  bold: true,\n  size: 14\n});\n\nawait context.sync();\n", "ps1": "# Append the \"Acknowledgments and Dedication\" heading to the last (empty)\n# paragraph of the document body, using the same Arial/bold/14pt styling\n# already used for the other title-page headings (e.g. \"Abstract\").\n$d = $word.ActiveDocument\n$p = $d.Paragraphs.Last\n$r = $p.Range\n$r.Text = \"Acknowledgments and Dedication\"\n$r.Font.Name = \"Arial\"\n$r.Font.NameAscii = \"Arial\"\n$r.Font.NameBi = \"Arial\"\n$r.Font.Bold = $true\n$r.Font.BoldBi = $true\n$r.Font.Size = 14\n$r.Font.SizeBi = 14\n"}
